# Generate Report for Handoff
# The e03c2c70-2abe-4799-a2b8-854130d6eb21 file was re-handed-off, so its
# "Latest Handoff Datetime" on each language sheet advances, and the
# Overview sheet's "Latest HO Xliff Generate Date" (the max across
# languages) advances to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2017-02-17 07:36:37"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2017-02-17 07:36:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2017-02-17 07:36:37"
